# Update the sequence-diagram deck:
#  1) Refresh every "automatically updated" date/time placeholder
#     (slide master, all slide layouts, and the notes master) to 9/9/2019.
#  2) Rename the "call mapping" label (controller -> logic call) to
#     "Call methods" and nudge its box slightly to the right.

$p = $ppt.ActivePresentation

$newDate = "9/9/2019"

# --- Slide master date placeholder -----------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every slide layout's date placeholder ----------------------------
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Notes master date placeholder ------------------------------------
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- Slide 1: "call mapping" textbox -> "Call methods" ----------------
$slide = $p.Slides.Item(1)
$callBox = $slide.Shapes.Item("TextBox 122")
$callBox.TextFrame.TextRange.Text = "Call methods"
$callBox.Left = 591.6989763779527
$callBox.Top = 236.88535433070865
